# Regenerate merged AHB files
# - Rename the "_old" / "_new" header-label suffixes to "_FV2310" / "_FV2404"
# - Freeze the header row (row 1)
# - Turn the data range into a native Excel Table (ListObject) with AutoFilter

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AHB-Diff")

# ---------------------------------------------------------------------------
# 1. Rename header labels: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}
# K1 ("diff") stays as-is.

# ---------------------------------------------------------------------------
# 2. Freeze panes at row 1 (top-left cell of the scrolling pane is A2)
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Convert the used range into a native Table (ListObject) with AutoFilter
# ---------------------------------------------------------------------------
$tblRange = $ws.Range("A1:U63")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tblRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

Write-Host "Edit complete."
